$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.831.79"
$ws.Range("E2").Value = "  -4.31%  "

$ws.Range("D3").Value = "3.407.01"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.76"
$ws.Range("B5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = "  -6.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "637.98"
$ws.Range("B6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = "  -2.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.42"
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.400"
$ws.Range("B8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "  -5.96%  "

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.961"
$ws.Range("B10").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "  -6.29%  "

$ws.Range("D11").Value = "3.408.93"
$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.199"
$ws.Range("B12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "  -4.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.27"
$ws.Range("B13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E13").Value = "  -4.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.18"
$ws.Range("B14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = "  +1.86%  "

$ws.Range("D15").Value = "93.871.96"
$ws.Range("E15").Value = "  -4.00%  "

$ws.Range("D16").Value = "4.042.38"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000249"
$ws.Range("B17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "  -3.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.24"
$ws.Range("B18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = "  -10.13%  "

$ws.Range("D19").Value = "3.407.90"
$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.30"
$ws.Range("B20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "  -3.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.47"
$ws.Range("B21").Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "  +2.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.471"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "  -8.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "495.58"
$ws.Range("B23").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.22"
$ws.Range("B24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4122) | Out-Null
$ws.Range("E24").Value = "  -5.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000191"
$ws.Range("B25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "  -4.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.45"
$ws.Range("B26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "  -5.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.99"
$ws.Range("B27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "  -5.63%  "

$ws.Range("D28").Value = "3.591.70"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.81"
$ws.Range("B29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = "  -4.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.46"
$ws.Range("B30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.72"
$ws.Range("B32").Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4122) | Out-Null
$ws.Range("E32").Value = "  +4.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.135"
$ws.Range("B33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "  -5.28%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("B34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.176"
$ws.Range("B35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = "  -7.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "29.30"
$ws.Range("B36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4122) | Out-Null
$ws.Range("E36").Value = "  +2.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.543"
$ws.Range("B37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4122) | Out-Null
$ws.Range("E37").Value = "  -2.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "539.98"
$ws.Range("B38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4122) | Out-Null
$ws.Range("E38").Value = "  +3.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.64"
$ws.Range("B39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4122) | Out-Null
$ws.Range("E39").Value = "  -3.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("B40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4122) | Out-Null
$ws.Range("E40").Value = "  -2.89%  "

$ws.Range("E41").Value = "  -0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.149"
$ws.Range("B42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.899"
$ws.Range("B43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4122) | Out-Null
$ws.Range("E43").Value = "  +6.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "24.02"
$ws.Range("B44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("B45").Value = "MantraDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.70"
$ws.Range("B45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("E45").Value = "  +0.57%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.70"
$ws.Range("B46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4122) | Out-Null
$ws.Range("E46").Value = "  -1.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.58"
$ws.Range("B47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4122) | Out-Null
$ws.Range("E47").Value = "  +0.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.18"
$ws.Range("B48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$ws.Range("E48").Value = "  -3.06%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.14"
$ws.Range("B49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4122) | Out-Null
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0404"
$ws.Range("B50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4122) | Out-Null
$ws.Range("E50").Value = "  -4.84%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.28"
$ws.Range("B51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").Value = "  +2.34%  "
